# Auto-applied edit matching the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap paired rows (full-row swap, column A/id untouched) ---
$rowA = $ws.Range("B121:AC121").Value2
$rowB = $ws.Range("B122:AC122").Value2
$ws.Range("B121:AC121").Value2 = $rowB
$ws.Range("B122:AC122").Value2 = $rowA

$rowA = $ws.Range("B125:AC125").Value2
$rowB = $ws.Range("B126:AC126").Value2
$ws.Range("B125:AC125").Value2 = $rowB
$ws.Range("B126:AC126").Value2 = $rowA

$rowA = $ws.Range("B128:AC128").Value2
$rowB = $ws.Range("B129:AC129").Value2
$ws.Range("B128:AC128").Value2 = $rowB
$ws.Range("B129:AC129").Value2 = $rowA

$rowA = $ws.Range("B137:AC137").Value2
$rowB = $ws.Range("B138:AC138").Value2
$ws.Range("B137:AC137").Value2 = $rowB
$ws.Range("B138:AC138").Value2 = $rowA

$rowA = $ws.Range("B140:AC140").Value2
$rowB = $ws.Range("B142:AC142").Value2
$ws.Range("B140:AC140").Value2 = $rowB
$ws.Range("B142:AC142").Value2 = $rowA

$rowA = $ws.Range("B150:AC150").Value2
$rowB = $ws.Range("B151:AC151").Value2
$ws.Range("B150:AC150").Value2 = $rowB
$ws.Range("B151:AC151").Value2 = $rowA

$rowA = $ws.Range("B157:AC157").Value2
$rowB = $ws.Range("B158:AC158").Value2
$ws.Range("B157:AC157").Value2 = $rowB
$ws.Range("B158:AC158").Value2 = $rowA

$rowA = $ws.Range("B164:AC164").Value2
$rowB = $ws.Range("B165:AC165").Value2
$ws.Range("B164:AC164").Value2 = $rowB
$ws.Range("B165:AC165").Value2 = $rowA

$rowA = $ws.Range("B191:AC191").Value2
$rowB = $ws.Range("B192:AC192").Value2
$ws.Range("B191:AC191").Value2 = $rowB
$ws.Range("B192:AC192").Value2 = $rowA

$rowA = $ws.Range("B194:AC194").Value2
$rowB = $ws.Range("B195:AC195").Value2
$ws.Range("B194:AC194").Value2 = $rowB
$ws.Range("B195:AC195").Value2 = $rowA

$rowA = $ws.Range("B210:AC210").Value2
$rowB = $ws.Range("B211:AC211").Value2
$ws.Range("B210:AC210").Value2 = $rowB
$ws.Range("B211:AC211").Value2 = $rowA

$rowA = $ws.Range("B216:AC216").Value2
$rowB = $ws.Range("B217:AC217").Value2
$ws.Range("B216:AC216").Value2 = $rowB
$ws.Range("B217:AC217").Value2 = $rowA

# --- Step 2: Replace content of rows 234 and 235 (new fixtures replacing old) ---
$ws.Range("B234").Value2 = 6940788
$ws.Range("E234").Value2 = 45387.66666666666
$ws.Range("F234").Value2 = "Al Ittihad Jeddah"
$ws.Range("G234").Value2 = "Al Taawon Buraidah"
$ws.Range("K234").Value2 = 1.727
$ws.Range("L234").Value2 = 3.6
$ws.Range("M234").Value2 = 4.333
$ws.Range("N234").Value2 = 1.727
$ws.Range("O234").Value2 = 3.6
$ws.Range("P234").Value2 = 4.333
$ws.Range("Q234").Value2 = -0.75
$ws.Range("R234").Value2 = 1.975
$ws.Range("S234").Value2 = 1.825
$ws.Range("T234").Value2 = 2.75
$ws.Range("U234").Value2 = 1.9
$ws.Range("V234").Value2 = 1.9
$ws.Range("W234").Value2 = 0
$ws.Range("X234").Value2 = 0
$ws.Range("Y234").Value2 = 0
$ws.Range("Z234").Value2 = 0
$ws.Range("AA234").Value2 = 0

$ws.Range("B235").Value2 = 6941439
$ws.Range("E235").Value2 = 45387.66666666666
$ws.Range("F235").Value2 = "Al Wehda Mecca"
$ws.Range("G235").Value2 = "Al Ahli Jeddah"
$ws.Range("K235").Value2 = 5
$ws.Range("L235").Value2 = 4.333
$ws.Range("M235").Value2 = 1.5
$ws.Range("N235").Value2 = 5
$ws.Range("O235").Value2 = 4.333
$ws.Range("P235").Value2 = 1.5
$ws.Range("Q235").Value2 = 1
$ws.Range("R235").Value2 = 2
$ws.Range("S235").Value2 = 1.8
$ws.Range("T235").Value2 = 3
$ws.Range("U235").Value2 = 1.95
$ws.Range("V235").Value2 = 1.85
$ws.Range("W235").Value2 = 0
$ws.Range("X235").Value2 = 0
$ws.Range("Y235").Value2 = 0
$ws.Range("Z235").Value2 = 0
$ws.Range("AA235").Value2 = 0

# --- Step 3: Append new rows 236-242 with fresh fixtures ---
# Copy formatting (styles) from row 235 down through the new rows
$ws.Range("A235:AC235").Copy() | Out-Null
$ws.Range("A236:AC242").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A236").Value2 = 234
$ws.Range("B236").Value2 = 6941441
$ws.Range("E236").Value2 = 45387.66666666666
$ws.Range("F236").Value2 = "Damac FC"
$ws.Range("G236").Value2 = "AlNassr Riyadh"
$ws.Range("K236").Value2 = 4.333
$ws.Range("L236").Value2 = 4.75
$ws.Range("M236").Value2 = 1.55
$ws.Range("N236").Value2 = 5
$ws.Range("O236").Value2 = 5
$ws.Range("P236").Value2 = 1.45
$ws.Range("Q236").Value2 = 1.25
$ws.Range("R236").Value2 = 1.85
$ws.Range("S236").Value2 = 1.95
$ws.Range("T236").Value2 = 3.25
$ws.Range("U236").Value2 = 1.9
$ws.Range("V236").Value2 = 1.9
$ws.Range("W236").Value2 = 0
$ws.Range("X236").Value2 = 0
$ws.Range("Y236").Value2 = 0
$ws.Range("Z236").Value2 = 0
$ws.Range("AA236").Value2 = 0

$ws.Range("A237").Value2 = 235
$ws.Range("B237").Value2 = 6941442
$ws.Range("E237").Value2 = 45387.66666666666
$ws.Range("F237").Value2 = "Al Khaleej Saihat"
$ws.Range("G237").Value2 = "Al Hilal Riyadh"
$ws.Range("K237").Value2 = 9.5
$ws.Range("L237").Value2 = 6
$ws.Range("M237").Value2 = 1.222
$ws.Range("N237").Value2 = 9.5
$ws.Range("O237").Value2 = 6
$ws.Range("P237").Value2 = 1.222
$ws.Range("Q237").Value2 = 1.75
$ws.Range("R237").Value2 = 1.975
$ws.Range("S237").Value2 = 1.825
$ws.Range("T237").Value2 = 3.25
$ws.Range("U237").Value2 = 1.95
$ws.Range("V237").Value2 = 1.85
$ws.Range("W237").Value2 = 0
$ws.Range("X237").Value2 = 0
$ws.Range("Y237").Value2 = 0
$ws.Range("Z237").Value2 = 0
$ws.Range("AA237").Value2 = 0

$ws.Range("A238").Value2 = 236
$ws.Range("B238").Value2 = 6941444
$ws.Range("E238").Value2 = 45388.66666666666
$ws.Range("F238").Value2 = "Al Ittifaq Dammam"
$ws.Range("G238").Value2 = "Al Riyadh"
$ws.Range("K238").Value2 = 1.3
$ws.Range("L238").Value2 = 4.75
$ws.Range("M238").Value2 = 8.5
$ws.Range("N238").Value2 = 1.333
$ws.Range("O238").Value2 = 4.75
$ws.Range("P238").Value2 = 7
$ws.Range("Q238").Value2 = -1.5
$ws.Range("R238").Value2 = 2.025
$ws.Range("S238").Value2 = 1.775
$ws.Range("T238").Value2 = 2.75
$ws.Range("U238").Value2 = 1.9
$ws.Range("V238").Value2 = 1.9
$ws.Range("W238").Value2 = 0
$ws.Range("X238").Value2 = 0
$ws.Range("Y238").Value2 = 0
$ws.Range("Z238").Value2 = 0
$ws.Range("AA238").Value2 = 0

$ws.Range("A239").Value2 = 237
$ws.Range("B239").Value2 = 6940826
$ws.Range("E239").Value2 = 45388.66666666666
$ws.Range("F239").Value2 = "Al Taee"
$ws.Range("G239").Value2 = "Al Shabab Riyadh"
$ws.Range("K239").Value2 = 3.8
$ws.Range("L239").Value2 = 3.6
$ws.Range("M239").Value2 = 1.8
$ws.Range("N239").Value2 = 4
$ws.Range("O239").Value2 = 3.6
$ws.Range("P239").Value2 = 1.75
$ws.Range("Q239").Value2 = 0.75
$ws.Range("R239").Value2 = 1.8
$ws.Range("S239").Value2 = 2
$ws.Range("T239").Value2 = 2.75
$ws.Range("U239").Value2 = 1.95
$ws.Range("V239").Value2 = 1.85
$ws.Range("W239").Value2 = 0
$ws.Range("X239").Value2 = 0
$ws.Range("Y239").Value2 = 0
$ws.Range("Z239").Value2 = 0
$ws.Range("AA239").Value2 = 0

$ws.Range("A240").Value2 = 238
$ws.Range("B240").Value2 = 6941443
$ws.Range("E240").Value2 = 45389.66666666666
$ws.Range("F240").Value2 = "Abha"
$ws.Range("G240").Value2 = "Al Fateh SC"
$ws.Range("K240").Value2 = 3
$ws.Range("L240").Value2 = 3.25
$ws.Range("M240").Value2 = 2.2
$ws.Range("N240").Value2 = 3.8
$ws.Range("O240").Value2 = 3.4
$ws.Range("P240").Value2 = 1.8
$ws.Range("Q240").Value2 = 0.5
$ws.Range("R240").Value2 = 1.9
$ws.Range("S240").Value2 = 1.9
$ws.Range("T240").Value2 = 2.75
$ws.Range("U240").Value2 = 1.8
$ws.Range("V240").Value2 = 2
$ws.Range("W240").Value2 = 0
$ws.Range("X240").Value2 = 0
$ws.Range("Y240").Value2 = 0
$ws.Range("Z240").Value2 = 0
$ws.Range("AA240").Value2 = 0

$ws.Range("A241").Value2 = 239
$ws.Range("B241").Value2 = 6941440
$ws.Range("E241").Value2 = 45389.66666666666
$ws.Range("F241").Value2 = "Al Raed"
$ws.Range("G241").Value2 = "Al Hazm"
$ws.Range("K241").Value2 = 1.65
$ws.Range("L241").Value2 = 3.75
$ws.Range("M241").Value2 = 4.75
$ws.Range("N241").Value2 = 1.65
$ws.Range("O241").Value2 = 3.75
$ws.Range("P241").Value2 = 4.75
$ws.Range("Q241").Value2 = -0.75
$ws.Range("R241").Value2 = 1.825
$ws.Range("S241").Value2 = 1.975
$ws.Range("T241").Value2 = 2.75
$ws.Range("U241").Value2 = 1.95
$ws.Range("V241").Value2 = 1.85
$ws.Range("W241").Value2 = 0
$ws.Range("X241").Value2 = 0
$ws.Range("Y241").Value2 = 0
$ws.Range("Z241").Value2 = 0
$ws.Range("AA241").Value2 = 0

$ws.Range("A242").Value2 = 240
$ws.Range("B242").Value2 = 6941438
$ws.Range("E242").Value2 = 45389.66666666666
$ws.Range("F242").Value2 = "Al Fayha"
$ws.Range("G242").Value2 = "Al Akhdoud"
$ws.Range("K242").Value2 = 2.5
$ws.Range("L242").Value2 = 3.2
$ws.Range("M242").Value2 = 2.625
$ws.Range("N242").Value2 = 2.5
$ws.Range("O242").Value2 = 3.2
$ws.Range("P242").Value2 = 2.625
$ws.Range("Q242").Value2 = 0
$ws.Range("R242").Value2 = 1.85
$ws.Range("S242").Value2 = 1.95
$ws.Range("T242").Value2 = 2.5
$ws.Range("U242").Value2 = 1.95
$ws.Range("V242").Value2 = 1.85
$ws.Range("W242").Value2 = 0
$ws.Range("X242").Value2 = 0
$ws.Range("Y242").Value2 = 0
$ws.Range("Z242").Value2 = 0
$ws.Range("AA242").Value2 = 0

# --- Step 4: Update sheet dimension to reflect new extent ---
$ws.Range("A1:AC242").Select() | Out-Null